$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424029619327598
$ws.Range("D2").Value = 0.1908596320807447
$ws.Range("E2").Value = 0.1668666624628585
$ws.Range("F2").Value = 1.416187686152867
$ws.Range("G2").Value = 0.002454875377159045
$ws.Range("J2").Value = 0.1877081658228192
$ws.Range("K2").Value = 0.7717270595578327
$ws.Range("N2").Value = 1.503750771621405
$ws.Range("O2").Value = 3.395160965842877

$ws.Range("B3").Value = 0.132947005609239
$ws.Range("D3").Value = 0.1858357018008832
$ws.Range("E3").Value = 0.1628844645285739
$ws.Range("F3").Value = 1.413534176495347
$ws.Range("G3").Value = 0.002457706336129722
$ws.Range("J3").Value = 0.1835469407849928
$ws.Range("K3").Value = 0.6879461982244948
$ws.Range("N3").Value = 1.518727073013128
$ws.Range("O3").Value = 3.400040374742616

$ws.Range("B4").Value = 0.1272099032646565
$ws.Range("D4").Value = 0.1828283054996405
$ws.Range("E4").Value = 0.1605231579989059
$ws.Range("F4").Value = 1.412704763756821
$ws.Range("G4").Value = 0.002459538163753612
$ws.Range("J4").Value = 0.1810993989830862
$ws.Range("K4").Value = 0.636512058673901
$ws.Range("N4").Value = 1.528481357279112
$ws.Range("O4").Value = 3.405027370312382

$ws.Range("B5").Value = 0.1248894866874934
$ws.Range("D5").Value = 0.1816222991821093
$ws.Range("E5").Value = 0.159582010501591
$ws.Range("F5").Value = 1.412567895188431
$ws.Range("G5").Value = 0.002460308257359896
$ws.Range("J5").Value = 0.1801290319338449
$ws.Range("K5").Value = 0.6155552397836743
$ws.Range("N5").Value = 1.532596892087085
$ws.Range("O5").Value = 3.407560029755786

$ws.Range("B6").Value = 0.1245052453511732
$ws.Range("D6").Value = 0.1814232250943348
$ws.Range("E6").Value = 0.159427009409896
$ws.Range("F6").Value = 1.412557315530506
$ws.Range("G6").Value = 0.002460437558941542
$ws.Range("J6").Value = 0.1799695356209554
$ws.Range("K6").Value = 0.6120755898822381
$ws.Range("N6").Value = 1.533288766785972
$ws.Range("O6").Value = 3.40801079302986

$ws.Range("B7").Value = 0.1271785382285202
$ws.Range("D7").Value = 0.1828119616720869
$ws.Range("E7").Value = 0.1605103798552854
$ws.Range("F7").Value = 1.412702103566097
$ws.Range("G7").Value = 0.002459548453720816
$ws.Range("J7").Value = 0.1810862028420601
$ws.Range("K7").Value = 0.6362294140536449
$ws.Range("N7").Value = 1.528536291496447
$ws.Range("O7").Value = 3.4050595007933

$ws.Range("B8").Value = 0.1391283515493882
$ws.Range("D8").Value = 0.189111385817867
$ws.Range("E8").Value = 0.165476232925279
$ws.Range("F8").Value = 1.415106738033714
$ws.Range("G8").Value = 0.002455832106355369
$ws.Range("J8").Value = 0.1862510703903979
$ws.Range("K8").Value = 0.7428384572560276
$ws.Range("N8").Value = 1.508798672913215
$ws.Range("O8").Value = 3.396429996944477

$ws.Range("B9").Value = 0.1631021056845725
$ws.Range("D9").Value = 0.2020748494859532
$ws.Range("E9").Value = 0.1758780626646228
$ws.Range("F9").Value = 1.426171281000663
$ws.Range("G9").Value = 0.002449283853570005
$ws.Range("J9").Value = 0.1972327299992855
$ws.Range("K9").Value = 0.9519216557770278
$ws.Range("N9").Value = 1.474523216390303
$ws.Range("O9").Value = 3.39532042279842

$ws.Range("B10").Value = 0.1810381931966276
$ws.Range("D10").Value = 0.2119680826587
$ws.Range("E10").Value = 0.1839247148822167
$ws.Range("F10").Value = 1.438177911036902
$ws.Range("G10").Value = 0.002444919131246361
$ws.Range("J10").Value = 0.205823245178081
$ws.Range("K10").Value = 1.105513991186115
$ws.Range("N10").Value = 1.452036246664363
$ws.Range("O10").Value = 3.404172024040747

$ws.Range("B11").Value = 0.1892665011768315
$ws.Range("D11").Value = 0.2165483077114629
$ws.Range("E11").Value = 0.187673174197684
$ws.Range("F11").Value = 1.444483729290923
$ws.Range("G11").Value = 0.002443029446684976
$ws.Range("J11").Value = 0.2098452123204737
$ws.Range("K11").Value = 1.175376011687888
$ws.Range("N11").Value = 1.442390484504479
$ws.Range("O11").Value = 3.410304123779781

$ws.Range("B12").Value = 0.1923921276373761
$ws.Range("D12").Value = 0.2182941120397004
$ws.Range("E12").Value = 0.189105250112803
$ws.Range("F12").Value = 1.446993014037844
$ws.Range("G12").Value = 0.002442327582168071
$ws.Range("J12").Value = 0.2113846505606034
$ws.Range("K12").Value = 1.201828909811184
$ws.Range("N12").Value = 1.43882175836567
$ws.Range("O12").Value = 3.412929355018008

$ws.Range("B13").Value = 0.1917185376891268
$ws.Range("D13").Value = 0.2179176175515209
$ws.Range("E13").Value = 0.1887962664621767
$ws.Range("F13").Value = 1.446447194203301
$ws.Range("G13").Value = 0.002442478132225161
$ws.Range("J13").Value = 0.2110523753176352
$ws.Range("K13").Value = 1.196131927675026
$ws.Range("N13").Value = 1.439586615494761
$ws.Range("O13").Value = 3.412350476363144

$ws.Range("B14").Value = 0.1895234541393052
$ws.Range("D14").Value = 0.2166917087036353
$ws.Range("E14").Value = 0.1877907393260543
$ws.Range("F14").Value = 1.444687736166713
$ws.Range("G14").Value = 0.002442971429219372
$ws.Range("J14").Value = 0.2099715340233388
$ws.Range("K14").Value = 1.177552361135326
$ws.Range("N14").Value = 1.442095201566218
$ws.Range("O14").Value = 3.410514025641817

$ws.Range("B15").Value = 0.1881801654824784
$ws.Range("D15").Value = 0.2159422823669104
$ws.Range("E15").Value = 0.1871764667273936
$ws.Range("F15").Value = 1.443625829747759
$ws.Range("G15").Value = 0.002443275372604498
$ws.Range("J15").Value = 0.2093116241179587
$ws.Range("K15").Value = 1.166171500883195
$ws.Range("N15").Value = 1.443642710695897
$ws.Range("O15").Value = 3.409428635233894

$ws.Range("B16").Value = 0.1805018279468413
$ws.Range("D16").Value = 0.2116703507454787
$ws.Range("E16").Value = 0.1836815117169621
$ws.Range("F16").Value = 1.437782797560999
$ws.Range("G16").Value = 0.002445044549818958
$ws.Range("J16").Value = 0.2055626960703592
$ws.Range("K16").Value = 1.100948082492721
$ws.Range("N16").Value = 1.452678363379029
$ws.Range("O16").Value = 3.403813676241953

$ws.Range("B17").Value = 0.1758089772235394
$ws.Range("D17").Value = 0.2090700160709247
$ws.Range("E17").Value = 0.1815599810487925
$ws.Range("F17").Value = 1.434414476203415
$ws.Range("G17").Value = 0.002446154385786335
$ws.Range("J17").Value = 0.2032920731297168
$ws.Range("K17").Value = 1.060932751737425
$ws.Range("N17").Value = 1.458370923186401
$ws.Range("O17").Value = 3.400908606882524

$ws.Range("B18").Value = 0.1731162882724675
$ws.Range("D18").Value = 0.2075818838471974
$ws.Range("E18").Value = 0.1803480176911521
$ws.Range("F18").Value = 1.432556539638838
$ws.Range("G18").Value = 0.002446801759922115
$ws.Range("J18").Value = 0.2019968104354035
$ws.Range("K18").Value = 1.037916319331117
$ws.Range("N18").Value = 1.461700071870695
$ws.Range("O18").Value = 3.399435841354375

$ws.Range("B19").Value = 0.1722057147579363
$ws.Range("D19").Value = 0.2070793208340405
$ws.Range("E19").Value = 0.1799390918331198
$ws.Range("F19").Value = 1.431941115370336
$ws.Range("G19").Value = 0.00244702250159799
$ws.Range("J19").Value = 0.2015601012903261
$ws.Range("K19").Value = 1.030123271780582
$ws.Range("N19").Value = 1.462836700788351
$ws.Range("O19").Value = 3.398971209901731

$ws.Range("B20").Value = 0.1763078663297506
$ws.Range("D20").Value = 0.209346049408353
$ws.Range("E20").Value = 0.1817849644230023
$ws.Range("F20").Value = 1.434764818398932
$ws.Range("G20").Value = 0.002446035308140693
$ws.Range("J20").Value = 0.203532673318378
$ws.Range("K20").Value = 1.065192530284605
$ws.Range("N20").Value = 1.457759254644223
$ws.Range("O20").Value = 3.401197345603919

$ws.Range("B21").Value = 0.1901679404370924
$ws.Range("D21").Value = 0.2170514799558418
$ws.Range("E21").Value = 0.1880857449375597
$ws.Range("F21").Value = 1.445201236255357
$ws.Range("G21").Value = 0.002442826164240101
$ws.Range("J21").Value = 0.2102885580483616
$ws.Range("K21").Value = 1.183009703039147
$ws.Range("N21").Value = 1.441356091724366
$ws.Range("O21").Value = 3.411045205657416

$ws.Range("B22").Value = 0.199283011550321
$ws.Range("D22").Value = 0.2221536680298044
$ws.Range("E22").Value = 0.1922771855530812
$ws.Range("F22").Value = 1.45272972215129
$ws.Range("G22").Value = 0.002440808733309314
$ws.Range("J22").Value = 0.2147995488694647
$ws.Range("K22").Value = 1.25999588506852
$ws.Range("N22").Value = 1.43112479009487
$ws.Range("O22").Value = 3.419248447325373

$ws.Range("B23").Value = 0.194413000510437
$ws.Range("D23").Value = 0.2194245041621627
$ws.Range("E23").Value = 0.1900334206941281
$ws.Range("F23").Value = 1.448646857783913
$ws.Range("G23").Value = 0.002441878181870967
$ws.Range("J23").Value = 0.2123831983294195
$ws.Range("K23").Value = 1.218908600375244
$ws.Range("N23").Value = 1.436540679022563
$ws.Range("O23").Value = 3.414708409291649

$ws.Range("B24").Value = 0.176082302062909
$ws.Range("D24").Value = 0.2092212334400898
$ws.Range("E24").Value = 0.1816832253504472
$ws.Range("F24").Value = 1.434606183995541
$ws.Range("G24").Value = 0.002446089114172023
$ws.Range("J24").Value = 0.2034238663622148
$ws.Range("K24").Value = 1.063266718843977
$ws.Range("N24").Value = 1.45803561424286
$ws.Range("O24").Value = 3.401066191964361

$ws.Range("B25").Value = 0.1565594340400622
$ws.Range("D25").Value = 0.1985028700235461
$ws.Range("E25").Value = 0.1729930521166807
$ws.Range("F25").Value = 1.422497600953434
$ws.Range("G25").Value = 0.0024509766325083
$ws.Range("J25").Value = 0.1941703282371634
$ws.Range("K25").Value = 0.895360278988619
$ws.Range("N25").Value = 1.442095201566218
$ws.Range("O25").Value = 3.410514025641817
